$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metadata field headers (columns E and F)
$ws.Range("E1").Value = "metadatafield4"
$ws.Range("F1").Value = "metadatafield5"

# Enter the new values first (apostrophe forces text so "true"/"false" stay strings,
# not booleans), then apply number formats afterwards so the style table stays tidy.
$ws.Range("D2").Value = "'true"
$ws.Range("E2").Value = "'A Test"
$ws.Range("D3").Value = "'false"
$ws.Range("E3").Value = "'Another Test"

$ws.Range("C2").NumberFormat = "yyyy/mm/dd"
$ws.Range("C2").Formula = "=DATE(24,1,4)"
$ws.Range("C3").NumberFormat = "yyyy/mm/dd"
$ws.Range("C3").Formula = "=DATE(24,12,31)"

$ws.Range("D2").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("D3").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"

# Widen the first column and move the selection to reflect the edited cell
$ws.Columns.Item(1).ColumnWidth = 28.9
$ws.Range("E2").Select()
